$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose "Neomorphic" character pattern states are not defined in the
# publication -> replace with the "N;0=0;1=1;2=?" pattern string.
$neomorphicRows = @(2, 5, 10, 12, 14, 19, 21, 27, 32)
foreach ($r in $neomorphicRows) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.ClearFormats()
    $cell.Value = "N;0=0;1=1;2=?"
}

# Row 24's "Transformational" character pattern states are not defined in the
# publication -> replace with the "T;0=0;1=1;2=?" pattern string.
$cell24 = $ws.Cells.Item(24, 2)
$cell24.ClearFormats()
$cell24.Value = "T;0=0;1=1;2=?"

$ws.Range("B28").Select()
